# Journal de travail - rectification d'erreur concernant le scale de la camera
# This script reproduces the row/cell edits for Tableau1 (rows 63-66 on
# sheet "Feuil1") and the associated shared-string text changes.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Feuil1" - contains Tableau1
$ws2 = $wb.Worksheets.Item(2)   # "Sheet1" - statistics (recalculated via formulas)

# ---------------------------------------------------------------------
# Row 63: add "Fin" (D63) and switch the subject from Documentation to
# Implémentation; update the description text in place.
# ---------------------------------------------------------------------
$ws1.Range("C63").Copy()
$ws1.Range("D63").PasteSpecial(-4122)      # xlPasteFormats - reuse the time style (s=9)
$ws1.Range("D63").Value = 0.51041666666666663

$ws1.Range("F63").Value = "Implémentation"

$ws1.Range("G63").Value = "Implémentation de la génération du mouvement des pièces"

# ---------------------------------------------------------------------
# Row 64: the start time ("Début") moved later in the day and the "Fin"
# value is now filled in; description text updated in place.
# ---------------------------------------------------------------------
$ws1.Range("C64").Value = 0.5625

$ws1.Range("C64").Copy()
$ws1.Range("D64").PasteSpecial(-4122)      # xlPasteFormats - reuse the time style (s=9)
$ws1.Range("D64").Value = 0.62847222222222221

$ws1.Range("G64").Value = "Implémentation de la saisie de la souris"

# ---------------------------------------------------------------------
# Row 65: brand-new journal entry (Date / Début / Fin / Sujet /
# Description) added below the existing ones.
# ---------------------------------------------------------------------
$ws1.Range("B64").Copy()
$ws1.Range("B65").PasteSpecial(-4122)      # xlPasteFormats - reuse the date style (s=8)
$ws1.Range("B65").Value = 45429

$ws1.Range("C64").Copy()
$ws1.Range("C65").PasteSpecial(-4122)      # xlPasteFormats - reuse the time style (s=9)
$ws1.Range("C65").Value = 0.63888888888888895

$ws1.Range("D64").Copy()
$ws1.Range("D65").PasteSpecial(-4122)      # xlPasteFormats - reuse the time style (s=9)
$ws1.Range("D65").Value = 0.70486111111111116

$ws1.Range("F64").Copy()
$ws1.Range("F65").PasteSpecial(-4122)      # xlPasteFormats - reuse F column style (s=1)
$ws1.Range("F65").Value = "Implémentation"

$ws1.Range("G64").Copy()
$ws1.Range("G65").PasteSpecial(-4122)      # xlPasteFormats - reuse G column style (s=2)
$ws1.Range("G65").Value = "Résolution de problème avec la saisie de la souris"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Row 66: this entry is no longer part of the journal - clear the
# leftover "Durée" formula so the row is blank again.
# ---------------------------------------------------------------------
$ws1.Range("E66").ClearContents()

# ---------------------------------------------------------------------
# Restore the view to the bottom of the table (scroll position + the
# active selection moved one row down together with the new entry).
# ---------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("F66").Select()

# Sheet2/"Sheet1" statistics (SUMIF totals, percentages) and the chart
# cache all derive from Tableau1 via formulas, so they recompute on
# their own once the workbook recalculates.
$excel.CalculateFull()
